$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.951.17"
$ws.Range("E2").Value = "  +1.16%  "
$ws.Range("D3").Value = "2.541.26"
$ws.Range("E3").Value = "  +0.54%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "592.32"
$ws.Range("E5").Value = "  +0.30%  "
$ws.Range("D6").Value = "173.60"
$ws.Range("E6").Value = "  +0.35%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D9").Value = "2.541.66"
$ws.Range("E9").Value = "  +0.60%  "
$ws.Range("E10").Value = "  +0.62%  "
$ws.Range("E11").Value = "  +1.52%  "
$ws.Range("D12").Value = "5.09"
$ws.Range("E12").Value = "  -1.16%  "
$ws.Range("D13").Value = "0.343"
$ws.Range("E13").Value = "  -0.02%  "
$ws.Range("B14").Value = "Binance-PegBSC-USD"
$ws.Range("C14").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D14").Value = "2.55"
$ws.Range("E14").Value = "  +155.48%  "
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").Value = "26.51"
$ws.Range("E15").Value = "  -0.28%  "
$ws.Range("E16").Value = "  -0.46%  "
$ws.Range("E17").Value = "  +0.95%  "
$ws.Range("D18").Value = "67.713.18"
$ws.Range("E18").Value = "  +0.97%  "
$ws.Range("D19").Value = "2.544.01"
$ws.Range("E19").Value = "  +1.02%  "
$ws.Range("E20").Value = "  +3.33%  "
$ws.Range("D21").Value = "7.96"
$ws.Range("E21").Value = "  -1.46%  "
$ws.Range("D22").Value = "369.27"
$ws.Range("E22").Value = "  +4.24%  "
$ws.Range("E23").Value = "  -0.60%  "
$ws.Range("D24").Value = "4.59"
$ws.Range("E24").Value = "  -0.81%  "
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").Value = "1.01"
$ws.Range("E25").Value = "  +1.15%  "
$ws.Range("B26").Value = "Litecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D26").Value = "71.60"
$ws.Range("E26").Value = "  +2.71%  "
$ws.Range("E27").Value = "  -3.38%  "
$ws.Range("D28").Value = "9.97"
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("D29").Value = "2.609.23"
$ws.Range("E29").Value = "  -1.68%  "
$ws.Range("D30").Value = "0.0₃0971"
$ws.Range("E30").Value = "  -0.72%  "
$ws.Range("D31").Value = "8.44"
$ws.Range("E31").Value = "  +3.61%  "
$ws.Range("D32").Value = "541.34"
$ws.Range("E32").Value = "  +1.63%  "
$ws.Range("E33").Value = "  -0.40%  "
$ws.Range("E34").Value = "  +1.14%  "
$ws.Range("E35").Value = "  -1.13%  "
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.06%  "
$ws.Range("D37").Value = "159.34"
$ws.Range("E37").Value = "  +1.10%  "
$ws.Range("E38").Value = "  -1.83%  "
$ws.Range("D39").Value = "19.18"
$ws.Range("E39").Value = "  +2.88%  "
$ws.Range("D40").Value = "18.62"
$ws.Range("E40").Value = "  +0.92%  "
$ws.Range("B41").Value = "PolygonEcosystemToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D41").Value = "0.352"
$ws.Range("E41").Value = "  -0.75%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D42").Value = "5.17"
$ws.Range("E42").Value = "  +0.71%  "
$ws.Range("D43").Value = "1.78"
$ws.Range("E43").Value = "  -0.48%  "
$ws.Range("D44").Value = "2.57"
$ws.Range("E44").Value = "  +3.03%  "
$ws.Range("E45").Value = "  +0.15%  "
$ws.Range("D46").Value = "39.26"
$ws.Range("D47").Value = "0.0₆0289"
$ws.Range("E47").Value = "  +3.96%  "
$ws.Range("D48").Value = "147.66"
$ws.Range("E48").Value = "  -0.94%  "
$ws.Range("D49").Value = "3.71"
$ws.Range("E49").Value = "  +0.47%  "
$ws.Range("D50").Value = "0.553"
$ws.Range("E50").Value = "  -0.59%  "
$ws.Range("E51").Value = "  +1.26%  "
